$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("E2").Value = 3
$ws.Range("G2").Value = 4.979788333333333
$ws.Range("H2").Value = 14.939365
$ws.Range("I2").Value = 0.129176854764059
$ws.Range("J2").Value = 0.129176854764059
$ws.Range("K2").Value = 3
$ws.Range("M2").Value = 31.632955
$ws.Range("N2").Value = 94.898865
$ws.Range("O2").Value = 0.5000340016649593
$ws.Range("P2").Value = 0.5000340016649593
$ws.Range("Q2").Value = 157.5254202578583
$ws.Range("R2").Value = 1417.728782320725
$ws.Range("S2").Value = 0.06459281961016566
$ws.Range("T2").Value = 0.06459281961016568

# Row 3
$ws.Range("E3").Value = 3
$ws.Range("G3").Value = 4.979788333333333
$ws.Range("H3").Value = 14.939365
$ws.Range("I3").Value = 0.129176854764059
$ws.Range("J3").Value = 0.129176854764059
$ws.Range("K3").Value = 3
$ws.Range("M3").Value = 8.622273333333332
$ws.Range("N3").Value = 25.86682
$ws.Range("O3").Value = 0.1362955132808722
$ws.Range("P3").Value = 0.1362955132808722
$ws.Range("Q3").Value = 42.93709615214444
$ws.Range("R3").Value = 386.4338653693
$ws.Range("S3").Value = 0.0176062257240761
$ws.Range("T3").Value = 0.0176062257240761

# Row 4
$ws.Range("E4").Value = 3
$ws.Range("G4").Value = 4.979788333333333
$ws.Range("H4").Value = 14.939365
$ws.Range("I4").Value = 0.129176854764059
$ws.Range("J4").Value = 0.129176854764059
$ws.Range("K4").Value = 3
$ws.Range("M4").Value = 18.70876033333333
$ws.Range("N4").Value = 56.12628100000001
$ws.Range("O4").Value = 0.2957364019791172
$ws.Range("P4").Value = 0.2957364019791172
$ws.Range("Q4").Value = 93.16566643906279
$ws.Range("R4").Value = 838.4909979515651
$ws.Range("S4").Value = 0.03820229824690178
$ws.Range("T4").Value = 0.03820229824690179

# Row 5
$ws.Range("E5").Value = 3
$ws.Range("G5").Value = 4.979788333333333
$ws.Range("H5").Value = 14.939365
$ws.Range("I5").Value = 0.129176854764059
$ws.Range("J5").Value = 0.129176854764059
$ws.Range("K5").Value = 3
$ws.Range("M5").Value = 4.297619333333333
$ws.Range("N5").Value = 12.892858
$ws.Range("O5").Value = 0.06793408307505136
$ws.Range("P5").Value = 0.06793408307505136
$ws.Range("Q5").Value = 21.40123461724111
$ws.Range("R5").Value = 192.61111155517
$ws.Range("S5").Value = 0.008775511182915426
$ws.Range("T5").Value = 0.008775511182915427

# Row 6
$ws.Range("E6").Value = 3
$ws.Range("G6").Value = 18.019504
$ws.Range("H6").Value = 54.058512
$ws.Range("I6").Value = 0.467430078412646
$ws.Range("J6").Value = 0.4674300784126461
$ws.Range("K6").Value = 3
$ws.Range("M6").Value = 31.632955
$ws.Range("N6").Value = 94.898865
$ws.Range("O6").Value = 0.5000340016649593
$ws.Range("P6").Value = 0.5000340016649593
$ws.Range("Q6").Value = 570.01015915432
$ws.Range("R6").Value = 5130.09143238888
$ws.Range("S6").Value = 0.2337309326072411
$ws.Range("T6").Value = 0.2337309326072411

# Row 7
$ws.Range("E7").Value = 3
$ws.Range("G7").Value = 18.019504
$ws.Range("H7").Value = 54.058512
$ws.Range("I7").Value = 0.467430078412646
$ws.Range("J7").Value = 0.4674300784126461
$ws.Range("K7").Value = 3
$ws.Range("M7").Value = 8.622273333333332
$ws.Range("N7").Value = 25.86682
$ws.Range("O7").Value = 0.1362955132808722
$ws.Range("P7").Value = 0.1362955132808722
$ws.Range("Q7").Value = 155.3690888190933
$ws.Range("R7").Value = 1398.32179937184
$ws.Range("S7").Value = 0.06370862246016994
$ws.Range("T7").Value = 0.06370862246016995

# Row 8
$ws.Range("E8").Value = 3
$ws.Range("G8").Value = 18.019504
$ws.Range("H8").Value = 54.058512
$ws.Range("I8").Value = 0.467430078412646
$ws.Range("J8").Value = 0.4674300784126461
$ws.Range("K8").Value = 3
$ws.Range("M8").Value = 18.70876033333333
$ws.Range("N8").Value = 56.12628100000001
$ws.Range("O8").Value = 0.2957364019791172
$ws.Range("P8").Value = 0.2957364019791172
$ws.Range("Q8").Value = 337.1225816615413
$ws.Range("R8").Value = 3034.103234953872
$ws.Range("S8").Value = 0.1382360895665725
$ws.Range("T8").Value = 0.1382360895665726

# Row 9
$ws.Range("E9").Value = 3
$ws.Range("G9").Value = 18.019504
$ws.Range("H9").Value = 54.058512
$ws.Range("I9").Value = 0.467430078412646
$ws.Range("J9").Value = 0.4674300784126461
$ws.Range("K9").Value = 3
$ws.Range("M9").Value = 4.297619333333333
$ws.Range("N9").Value = 12.892858
$ws.Range("O9").Value = 0.06793408307505136
$ws.Range("P9").Value = 0.06793408307505136
$ws.Range("Q9").Value = 77.44096876747734
$ws.Range("R9").Value = 696.9687189072961
$ws.Range("S9").Value = 0.03175443377866247
$ws.Range("T9").Value = 0.03175443377866247

# Row 10
$ws.Range("E10").Value = 3
$ws.Range("G10").Value = 8.752692000000001
$ws.Range("H10").Value = 26.258076
$ws.Range("I10").Value = 0.2270468436801446
$ws.Range("J10").Value = 0.2270468436801446
$ws.Range("K10").Value = 3
$ws.Range("M10").Value = 31.632955
$ws.Range("N10").Value = 94.898865
$ws.Range("O10").Value = 0.5000340016649593
$ws.Range("P10").Value = 0.5000340016649593
$ws.Range("Q10").Value = 276.87351216486
$ws.Range("R10").Value = 2491.86160948374
$ws.Range("S10").Value = 0.1135311418107812
$ws.Range("T10").Value = 0.1135311418107812

# Row 11
$ws.Range("E11").Value = 3
$ws.Range("G11").Value = 8.752692000000001
$ws.Range("H11").Value = 26.258076
$ws.Range("I11").Value = 0.2270468436801446
$ws.Range("J11").Value = 0.2270468436801446
$ws.Range("K11").Value = 3
$ws.Range("M11").Value = 8.622273333333332
$ws.Range("N11").Value = 25.86682
$ws.Range("O11").Value = 0.1362955132808722
$ws.Range("P11").Value = 0.1362955132808722
$ws.Range("Q11").Value = 75.46810282648
$ws.Range("R11").Value = 679.21292543832
$ws.Range("S11").Value = 0.03094546609818727
$ws.Range("T11").Value = 0.03094546609818727

# Row 12
$ws.Range("E12").Value = 3
$ws.Range("G12").Value = 8.752692000000001
$ws.Range("H12").Value = 26.258076
$ws.Range("I12").Value = 0.2270468436801446
$ws.Range("J12").Value = 0.2270468436801446
$ws.Range("K12").Value = 3
$ws.Range("M12").Value = 18.70876033333333
$ws.Range("N12").Value = 56.12628100000001
$ws.Range("O12").Value = 0.2957364019791172
$ws.Range("P12").Value = 0.2957364019791172
$ws.Range("Q12").Value = 163.752016899484
$ws.Range("R12").Value = 1473.768152095356
$ws.Range("S12").Value = 0.06714601663068101
$ws.Range("T12").Value = 0.06714601663068102

# Row 13
$ws.Range("E13").Value = 3
$ws.Range("G13").Value = 8.752692000000001
$ws.Range("H13").Value = 26.258076
$ws.Range("I13").Value = 0.2270468436801446
$ws.Range("J13").Value = 0.2270468436801446
$ws.Range("K13").Value = 3
$ws.Range("M13").Value = 4.297619333333333
$ws.Range("N13").Value = 12.892858
$ws.Range("O13").Value = 0.06793408307505136
$ws.Range("P13").Value = 0.06793408307505136
$ws.Range("Q13").Value = 37.61573835791201
$ws.Range("R13").Value = 338.5416452212081
$ws.Range("S13").Value = 0.01542421914049514
$ws.Range("T13").Value = 0.01542421914049514

# Row 14
$ws.Range("E14").Value = 3
$ws.Range("G14").Value = 6.798175000000001
$ws.Range("H14").Value = 20.394525
$ws.Range("I14").Value = 0.1763462231431503
$ws.Range("J14").Value = 0.1763462231431503
$ws.Range("K14").Value = 3
$ws.Range("M14").Value = 31.632955
$ws.Range("N14").Value = 94.898865
$ws.Range("O14").Value = 0.5000340016649593
$ws.Range("P14").Value = 0.5000340016649593
$ws.Range("Q14").Value = 215.046363857125
$ws.Range("R14").Value = 1935.417274714125
$ws.Range("S14").Value = 0.08817910763677131
$ws.Range("T14").Value = 0.08817910763677132

# Row 15
$ws.Range("E15").Value = 3
$ws.Range("G15").Value = 6.798175000000001
$ws.Range("H15").Value = 20.394525
$ws.Range("I15").Value = 0.1763462231431503
$ws.Range("J15").Value = 0.1763462231431503
$ws.Range("K15").Value = 3
$ws.Range("M15").Value = 8.622273333333332
$ws.Range("N15").Value = 25.86682
$ws.Range("O15").Value = 0.1362955132808722
$ws.Range("P15").Value = 0.1362955132808722
$ws.Range("Q15").Value = 58.61572301783333
$ws.Range("R15").Value = 527.5415071605
$ws.Range("S15").Value = 0.0240351989984389
$ws.Range("T15").Value = 0.02403519899843891

# Row 16
$ws.Range("E16").Value = 3
$ws.Range("G16").Value = 6.798175000000001
$ws.Range("H16").Value = 20.394525
$ws.Range("I16").Value = 0.1763462231431503
$ws.Range("J16").Value = 0.1763462231431503
$ws.Range("K16").Value = 3
$ws.Range("M16").Value = 18.70876033333333
$ws.Range("N16").Value = 56.12628100000001
$ws.Range("O16").Value = 0.2957364019791172
$ws.Range("P16").Value = 0.2957364019791172
$ws.Range("Q16").Value = 127.1854267790583
$ws.Range("R16").Value = 1144.668841011525
$ws.Range("S16").Value = 0.0521519975349618
$ws.Range("T16").Value = 0.05215199753496182

# Row 17
$ws.Range("E17").Value = 3
$ws.Range("G17").Value = 6.798175000000001
$ws.Range("H17").Value = 20.394525
$ws.Range("I17").Value = 0.1763462231431503
$ws.Range("J17").Value = 0.1763462231431503
$ws.Range("K17").Value = 3
$ws.Range("M17").Value = 4.297619333333333
$ws.Range("N17").Value = 12.892858
$ws.Range("O17").Value = 0.06793408307505136
$ws.Range("P17").Value = 0.06793408307505136
$ws.Range("Q17").Value = 29.21596831138334
$ws.Range("R17").Value = 262.94371480245
$ws.Range("S17").Value = 0.01197991897297832
$ws.Range("T17").Value = 0.01197991897297832
